$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.Iteration = $true
$excel.MaxChange = 0.0001

$ws.Range("A4").Value = "FIC"
$ws.Range("A5").Value = "LBD"
